# "Chore: all excercises complete"
# Fills in the remaining VLOOKUP / HLOOKUP / XLOOKUP / INDEX-MATCH practice
# formulas in the workbook (the exercises were left blank before), and
# leaves the selection/navigation state the way the author left it
# (finishing up on the INDEX worksheet).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# VLOOKUP sheet
# ---------------------------------------------------------------------
$wsV = $wb.Worksheets.Item("VLOOKUP")
$wsV.Activate()

$wsV.Range("G9").Formula  = '=VLOOKUP(G6,B5:E13,4,)'
$wsV.Range("J14").Formula = '=VLOOKUP(G6,B5:E13,2,FALSE)'

$wsV.Range("B5:E13").Select()

# ---------------------------------------------------------------------
# HLOOKUP sheet
# ---------------------------------------------------------------------
$wsH = $wb.Worksheets.Item("HLOOKUP")
$wsH.Activate()

$wsH.Range("L21").Formula = '=HLOOKUP(VLOOKUP(HLOOKUP!K6,VLOOKUP!B5:E13,3,FALSE),HLOOKUP!B5:I8,4)'

$wsH.Range("L22").Select()

# ---------------------------------------------------------------------
# XLOOKUP sheet
# ---------------------------------------------------------------------
$wsX = $wb.Worksheets.Item("XLOOKUP")
$wsX.Activate()

$wsX.Range("G9").Formula  = '=XLOOKUP(G6,E5:E13,D5:D13,"",)'
$wsX.Range("J14").Formula = '=XLOOKUP(G6,E5:E13,B5:B13,0,)'
$wsX.Range("L45").Formula = '=XLOOKUP(XLOOKUP(K30,E5:E13,C5:C13,0,0),C29:I29,C32:I32,0,-1)'

$excel.ActiveWindow.Zoom = 115
$wsX.Range("F38").Select()

# ---------------------------------------------------------------------
# INDEX sheet (INDEX/MATCH, plus an XLOOKUP-based alternative)
# ---------------------------------------------------------------------
$wsI = $wb.Worksheets.Item("INDEX")
$wsI.Activate()

$wsI.Range("K14").Formula = '=INDEX($C$6:$G$16,MATCH(J14,$C$6:$C$16,0), MATCH(I14,C6:$G$6,0))'
$wsI.Range("K15").FormulaArray = '=XLOOKUP(J15, C6:C16,XLOOKUP(I15,C6:G6,C6:G16))'
$wsI.Range("K16").FormulaArray = '=XLOOKUP(J16,C6:C16,XLOOKUP(I16,C6:G6,C6:G16))'

# Second approach shown side-by-side for comparison
$wsI.Range("M16").Value   = "Index match alt"
$wsI.Range("N16").Formula = '=INDEX($C$6:$G$16,MATCH(J16,$C$6:$C$16,0), MATCH(I16,C6:$G$6,0))'
$wsI.Columns.Item(13).AutoFit() | Out-Null

$wsI.Range("J25").Select()
